$wb = $excel.ActiveWorkbook

# --- Sheet: P_valores ---
$wsP = $wb.Worksheets.Item("P_valores")

$wsP.Range("C2").Value = 0.8473667966228036
$wsP.Range("D2").Value = 0.9341805931064964
$wsP.Range("E2").Value = 0.3057721417331292
$wsP.Range("F2").Value = 0.233290907269317

$wsP.Range("B3").Value = 0.8473667966228036
$wsP.Range("D3").Value = 0.9134837200885222
$wsP.Range("E3").Value = 0.5625568701275099
$wsP.Range("F3").Value = 0.5261978668931748

$wsP.Range("B4").Value = 0.9341805931064964
$wsP.Range("C4").Value = 0.9134837200885222
$wsP.Range("E4").Value = 0.175831384098754
$wsP.Range("F4").Value = 0.5266221328077765

$wsP.Range("B5").Value = 0.3057721417331292
$wsP.Range("C5").Value = 0.5625568701275099
$wsP.Range("D5").Value = 0.175831384098754
$wsP.Range("F5").Value = 0.9176839217878889

$wsP.Range("B6").Value = 0.233290907269317
$wsP.Range("C6").Value = 0.5261978668931748
$wsP.Range("D6").Value = 0.5266221328077765
$wsP.Range("E6").Value = 0.9176839217878889

# --- Sheet: Estadisticos_DM ---
$wsD = $wb.Worksheets.Item("Estadisticos_DM")

$wsD.Range("C2").Value = -0.1960779145433331
$wsD.Range("D2").Value = -0.08408345448158552
$wsD.Range("E2").Value = -1.063005742191642
$wsD.Range("F2").Value = -1.245780504086022

$wsD.Range("B3").Value = 0.1960779145433331
$wsD.Range("D3").Value = 0.1106253065631983
$wsD.Range("E3").Value = -0.5931256949115136
$wsD.Range("F3").Value = -0.6500273470049559

$wsD.Range("B4").Value = 0.08408345448158552
$wsD.Range("C4").Value = -0.1106253065631983
$wsD.Range("E4").Value = -1.425835432289283
$wsD.Range("F4").Value = -0.6493509535679736

$wsD.Range("B5").Value = 1.063005742191642
$wsD.Range("C5").Value = 0.5931256949115136
$wsD.Range("D5").Value = 1.425835432289283
$wsD.Range("F5").Value = 0.1052328286658491

$wsD.Range("B6").Value = 1.245780504086022
$wsD.Range("C6").Value = 0.6500273470049559
$wsD.Range("D6").Value = 0.6493509535679736
$wsD.Range("E6").Value = -0.1052328286658491

$wb.Save()
